$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 464, pushing the existing data (464:496) down to (467:499).
$ws.Rows("464:466").Insert()

# Row 464: Pintón, Volumen 80, Precio 14000 / 14000 / 14000, Precio $/Kg 700
$ws.Range("A464").Value = 8
$ws.Range("B464").Value = "Terminal La Palmera de La Serena"
$ws.Range("C464").Value = "Coquimbo"
$ws.Range("D464").Value = 44585
$ws.Range("E464").Value = 4
$ws.Range("F464").Value = "Fruta"
$ws.Range("G464").Value = 100108
$ws.Range("H464").Value = "Tropicales y subtropicales"
$ws.Range("I464").Value = 100108006
$ws.Range("J464").Value = "Plátano"
$ws.Range("K464").Value = "Sin especificar"
$ws.Range("L464").Value = "Pintón"
$ws.Range("M464").Value = 80
$ws.Range("N464").Value = 14000
$ws.Range("O464").Value = 14000
$ws.Range("P464").Value = 14000
$ws.Range("Q464").Value = "$/caja 20 kilos"
$ws.Range("R464").Value = "Ecuador"
$ws.Range("S464").Value = 700
$ws.Range("T464").Value = 20

# Row 465: Primera Maduro, Volumen 120, Precio 16000 / 16000 / 16000, Precio $/Kg 800
$ws.Range("A465").Value = 8
$ws.Range("B465").Value = "Terminal La Palmera de La Serena"
$ws.Range("C465").Value = "Coquimbo"
$ws.Range("D465").Value = 44585
$ws.Range("E465").Value = 4
$ws.Range("F465").Value = "Fruta"
$ws.Range("G465").Value = 100108
$ws.Range("H465").Value = "Tropicales y subtropicales"
$ws.Range("I465").Value = 100108006
$ws.Range("J465").Value = "Plátano"
$ws.Range("K465").Value = "Sin especificar"
$ws.Range("L465").Value = "Primera Maduro"
$ws.Range("M465").Value = 120
$ws.Range("N465").Value = 16000
$ws.Range("O465").Value = 16000
$ws.Range("P465").Value = 16000
$ws.Range("Q465").Value = "$/caja 20 kilos"
$ws.Range("R465").Value = "Ecuador"
$ws.Range("S465").Value = 800
$ws.Range("T465").Value = 20

# Row 466: Primera Pintón, Volumen 120, Precio 17000 / 17000 / 17000, Precio $/Kg 850
$ws.Range("A466").Value = 8
$ws.Range("B466").Value = "Terminal La Palmera de La Serena"
$ws.Range("C466").Value = "Coquimbo"
$ws.Range("D466").Value = 44585
$ws.Range("E466").Value = 4
$ws.Range("F466").Value = "Fruta"
$ws.Range("G466").Value = 100108
$ws.Range("H466").Value = "Tropicales y subtropicales"
$ws.Range("I466").Value = 100108006
$ws.Range("J466").Value = "Plátano"
$ws.Range("K466").Value = "Sin especificar"
$ws.Range("L466").Value = "Primera Pintón"
$ws.Range("M466").Value = 120
$ws.Range("N466").Value = 17000
$ws.Range("O466").Value = 17000
$ws.Range("P466").Value = 17000
$ws.Range("Q466").Value = "$/caja 20 kilos"
$ws.Range("R466").Value = "Ecuador"
$ws.Range("S466").Value = 850
$ws.Range("T466").Value = 20

# Ensure the date column keeps the date number format used by the rest of column D.
$ws.Range("D464:D466").NumberFormat = $ws.Range("D467").NumberFormat()
